$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix requirement-id labels: "MJ-0X" -> "MJ0X"
$ws.Range("C2").Value = "MJ01"
$ws.Range("C3").Value = "MJ02"
$ws.Range("C4").Value = "MJ03"
$ws.Range("C5").Value = "MJ04"
$ws.Range("C6").Value = "MJ05"
$ws.Range("C7").Value = "MJ06"
$ws.Range("C8").Value = "MJ07"

# 2. Update the UT Mapping text in F4 and restyle it (smaller, black, Arial font)
$ws.Range("F4").Value = "UT_TEST_CASE(1-10)"
$ws.Range("F4").Font.Color = 0
$ws.Range("F4").Font.Size = 10

# 3. Re-apply the standard wrapped/bordered formatting used by rows 2-6 to C7/C8
#    (their border previously didn't match the rest of the column)
$ws.Range("C4").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4. Add the new MJ08 requirement row
$ws.Range("C4").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("C9").Value = "MJ08"

$ws.Range("D8").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("D9").Value = "3.1.8"

$ws.Range("D8").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("E9").Font.Size = 11
$ws.Range("E9").Value = "NA"

$ws.Range("D8").Copy()
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("F9").Font.Size = 11
$ws.Range("F9").Value = "NA"

$ws.Range("G8").Copy()
$ws.Range("G9").PasteSpecial(-4122)
$ws.Range("G9").Value = "NA"

$excel.CutCopyMode = 0

# 5. Leave the selection where the last edit happened
$ws.Range("E9").Select()
